# Apply the "6.0.0" release update to the StructureDefinition-match-criteria
# workbook:
#   - Metadata sheet: bump Version, Date, set a real Publisher, replace the
#     bogus duplicated "Contact" rows with a single "Jurisdiction" row, and
#     remove the now-redundant duplicate row entirely.
#   - Elements sheet: give the root Extension row a real Short/Definition
#     (was showing the generic "Extension" / "An Extension" placeholder).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Row 11 is an exact duplicate of row 10 ("Contact" / "No display for
# ContactDetail"). Delete it outright - this shifts rows 12-21 up by one
# and shrinks the used range from A1:B21 to A1:B20.
$meta.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Date: refresh the publication timestamp
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher: was blank, now populated
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# Former duplicate "Contact" row (row 10) becomes the Jurisdiction row
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# Elements sheet: root Extension row (row 2) gets real Short/Definition
# text instead of the generic placeholders.
$elements = $wb.Worksheets.Item("Elements")
$elements.Cells.Item(2, 11).Value = "Match Criteria"
$elements.Cells.Item(2, 12).Value = "Criteria or algorythm used to identify the matched resource"
